# Added Boats Search test case on API level
#
# Summary of the change (per the target diff):
#  - Sheet "API" gains two new columns ("Country Name" / "City Name") right
#    after "TC ID/Name", reusing the old "Boat Name"/"Boat Type" cells'
#    positions & styles with new text ("turkey" / "istanbul").
#  - A brand-new "Boat Name" column is inserted after "City Name" holding a
#    new expected value "Motor-boat-Custom-Built---2005-refit-2017-".
#  - The "API" sheet becomes the active/selected sheet & tab, zoomed to 160%,
#    with D3 selected.
#  - The "GUI" sheet reverts to not being the selected tab, with the full
#    A1:XFD2 range selected (its data/content is unchanged).

$wb = $excel.ActiveWorkbook

$wsApi = $wb.Worksheets.Item("API")
$wsGui = $wb.Worksheets.Item("GUI")

# --- 1) "API" sheet: rewrite the header row & insert the new "Boat Name" column ---

# B/C currently hold "Boat Name"/"Boat Type" (header) and the matching
# values below them. Re-purpose them in place (keeps their existing
# style) as the new "Country Name"/"City Name" columns.
$wsApi.Range("B1").Value = "Country Name"
$wsApi.Range("C1").Value = "City Name"
$wsApi.Range("B2").Value = "turkey"
$wsApi.Range("C2").Value = "istanbul"

# Insert a fresh column in position D for the new "Boat Name" test data.
$wsApi.Range("D1:D1").EntireColumn.Insert() | Out-Null
$wsApi.Range("D1").Value = "Boat Name"
$wsApi.Range("D2").Value = "Motor-boat-Custom-Built---2005-refit-2017-"
# The newly inserted data cell carries no explicit style (unlike its
# neighbours), matching the default "Normal" style.
$wsApi.Range("D2").Style = "Normal"

# Best-effort column widths for the columns whose contents changed, so the
# sheet does not keep the old (now much too wide/narrow) best-fit widths.
$wsApi.Range("B1").ColumnWidth = 12.081333
$wsApi.Range("C1").ColumnWidth = 8.414667
$wsApi.Range("D1").ColumnWidth = 36.414667

# --- 2) View/selection housekeeping ---

# Make "API" the active sheet/tab, set its zoom and selection.
$wsApi.Activate()
$excel.ActiveWindow.Zoom = 160
$wsApi.Range("D3").Select() | Out-Null

# "GUI" goes back to a plain full-range selection and loses the active tab.
$wsGui.Range("A1:XFD2").Select() | Out-Null
$wsApi.Activate()

Write-Host "Boats Search API test case added"
